# Weekly update: insert a new week's price data for
# "Hortaliza, Terminal La Palmera de La Serena - Pepino dulce" before the
# previous week's rows (which shift down).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank rows starting at row 348; this pushes the existing
# rows 348-350 down to 351-353, and copies formatting (incl. the date style)
# from the row above.
$ws.Rows.Item(348).Insert()
$ws.Rows.Item(348).Insert()
$ws.Rows.Item(348).Insert()

# New week's data (fecha 44656) for the three quality grades.
$newRows = @(
    @{ Row = 348; Calidad = "Primera"; Volumen = 400; Min = 9000; Max = 10000; Prom = 9500; PKg = 528 },
    @{ Row = 349; Calidad = "Segunda"; Volumen = 300; Min = 7000; Max = 8000;  Prom = 7500; PKg = 417 },
    @{ Row = 350; Calidad = "Tercera"; Volumen = 240; Min = 5000; Max = 6000;  Prom = 5500; PKg = 306 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2 = 8
    $ws.Cells.Item($row, 2).Value2 = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($row, 3).Value2 = "Coquimbo"
    $ws.Cells.Item($row, 4).Value2 = 44656
    $ws.Cells.Item($row, 5).Value2 = 4
    $ws.Cells.Item($row, 6).Value2 = 100112043
    $ws.Cells.Item($row, 7).Value2 = "Pepino dulce"
    $ws.Cells.Item($row, 8).Value2 = "Cultivar IV Región"
    $ws.Cells.Item($row, 9).Value2 = $r.Calidad
    $ws.Cells.Item($row, 10).Value2 = $r.Volumen
    $ws.Cells.Item($row, 11).Value2 = $r.Min
    $ws.Cells.Item($row, 12).Value2 = $r.Max
    $ws.Cells.Item($row, 13).Value2 = $r.Prom
    $ws.Cells.Item($row, 14).Value2 = "`$/bandeja 18 kilos"
    $ws.Cells.Item($row, 15).Value2 = "Provincia de Limarí"
    $ws.Cells.Item($row, 16).Value2 = $r.PKg
    $ws.Cells.Item($row, 17).Value2 = 18
    $ws.Cells.Item($row, 18).Value2 = "Hortaliza"
}
